$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update constant columns B (time in ms) and C (auto capacity) for all existing data rows (2-51) ---
$ws.Range("B2:B51").Value = 33.94444444444444
$ws.Range("C2:C51").Value = 1.95

# --- Per-row D (auto scs / keys) and E (probability values) corrections for rows 2-51 ---
$ws.Cells.Item(2, 5).Value = 0.148
$ws.Cells.Item(4, 5).Value = 0.003
$ws.Cells.Item(5, 5).Value = 0.011
$ws.Cells.Item(6, 5).Value = 0.025
$ws.Cells.Item(7, 5).Value = 0.036
$ws.Cells.Item(8, 5).Value = 0.052
$ws.Cells.Item(9, 5).Value = 0.043
$ws.Cells.Item(11, 5).Value = 0.039
$ws.Cells.Item(12, 5).Value = 0.028
$ws.Cells.Item(13, 5).Value = 0.027
$ws.Cells.Item(14, 5).Value = 0.033
$ws.Cells.Item(15, 5).Value = 0.033
$ws.Cells.Item(16, 5).Value = 0.041
$ws.Cells.Item(17, 5).Value = 0.032
$ws.Cells.Item(18, 5).Value = 0.04
$ws.Cells.Item(19, 5).Value = 0.028
$ws.Cells.Item(20, 5).Value = 0.024
$ws.Cells.Item(21, 5).Value = 0.022
$ws.Cells.Item(22, 5).Value = 0.026
$ws.Cells.Item(23, 5).Value = 0.021
$ws.Cells.Item(24, 5).Value = 0.013
$ws.Cells.Item(25, 5).Value = 0.02
$ws.Cells.Item(26, 5).Value = 0.02
$ws.Cells.Item(27, 5).Value = 0.016
$ws.Cells.Item(29, 5).Value = 0.018
$ws.Cells.Item(30, 5).Value = 0.017
$ws.Cells.Item(31, 5).Value = 0.015
$ws.Cells.Item(32, 5).Value = 0.013
$ws.Cells.Item(33, 5).Value = 0.011
$ws.Cells.Item(34, 5).Value = 0.006
$ws.Cells.Item(35, 5).Value = 0.012
$ws.Cells.Item(36, 5).Value = 0.005
$ws.Cells.Item(38, 5).Value = 0.007
$ws.Cells.Item(39, 5).Value = 0.003
$ws.Cells.Item(40, 5).Value = 0.002
$ws.Cells.Item(41, 5).Value = 0.005
$ws.Cells.Item(42, 5).Value = 0.004
$ws.Cells.Item(43, 5).Value = 0.002
$ws.Cells.Item(44, 4).Value = 43
$ws.Cells.Item(44, 5).Value = 0.003
$ws.Cells.Item(45, 4).Value = 44
$ws.Cells.Item(45, 5).Value = 0.001
$ws.Cells.Item(46, 4).Value = 45
$ws.Cells.Item(46, 5).Value = 0.001
$ws.Cells.Item(47, 4).Value = 47
$ws.Cells.Item(47, 5).Value = 0.004
$ws.Cells.Item(48, 4).Value = 48
$ws.Cells.Item(48, 5).Value = 0.002
$ws.Cells.Item(49, 4).Value = 49
$ws.Cells.Item(49, 5).Value = 0.004
$ws.Cells.Item(50, 4).Value = 50
$ws.Cells.Item(50, 5).Value = 0.002
$ws.Cells.Item(51, 4).Value = 52
$ws.Cells.Item(51, 5).Value = 0.002

# --- Append new rows 52-59 (extend table; copy row-51 formatting for column A) ---
$ws.Cells.Item(52, 1).Value = 50
$ws.Cells.Item(52, 2).Value = 33.94444444444444
$ws.Cells.Item(52, 3).Value = 1.95
$ws.Cells.Item(52, 4).Value = 53
$ws.Cells.Item(52, 5).Value = 0.002
$ws.Cells.Item(53, 1).Value = 51
$ws.Cells.Item(53, 2).Value = 33.94444444444444
$ws.Cells.Item(53, 3).Value = 1.95
$ws.Cells.Item(53, 4).Value = 54
$ws.Cells.Item(53, 5).Value = 0.002
$ws.Cells.Item(54, 1).Value = 52
$ws.Cells.Item(54, 2).Value = 33.94444444444444
$ws.Cells.Item(54, 3).Value = 1.95
$ws.Cells.Item(54, 4).Value = 55
$ws.Cells.Item(54, 5).Value = 0.001
$ws.Cells.Item(55, 1).Value = 53
$ws.Cells.Item(55, 2).Value = 33.94444444444444
$ws.Cells.Item(55, 3).Value = 1.95
$ws.Cells.Item(55, 4).Value = 57
$ws.Cells.Item(55, 5).Value = 0.002
$ws.Cells.Item(56, 1).Value = 54
$ws.Cells.Item(56, 2).Value = 33.94444444444444
$ws.Cells.Item(56, 3).Value = 1.95
$ws.Cells.Item(56, 4).Value = 60
$ws.Cells.Item(56, 5).Value = 0.001
$ws.Cells.Item(57, 1).Value = 55
$ws.Cells.Item(57, 2).Value = 33.94444444444444
$ws.Cells.Item(57, 3).Value = 1.95
$ws.Cells.Item(57, 4).Value = 61
$ws.Cells.Item(57, 5).Value = 0.001
$ws.Cells.Item(58, 1).Value = 56
$ws.Cells.Item(58, 2).Value = 33.94444444444444
$ws.Cells.Item(58, 3).Value = 1.95
$ws.Cells.Item(58, 4).Value = 62
$ws.Cells.Item(58, 5).Value = 0.001
$ws.Cells.Item(59, 1).Value = 57
$ws.Cells.Item(59, 2).Value = 33.94444444444444
$ws.Cells.Item(59, 3).Value = 1.95
$ws.Cells.Item(59, 4).Value = 63
$ws.Cells.Item(59, 5).Value = 0.001

# --- Copy column-A cell formatting (bold, border, centered) from row 51 down through new rows 52-59 ---
$ws.Range("A51").Copy()
$ws.Range("A52:A59").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Output "edit complete"
